# Auto-generated edit script applying the Typhon_Profits profit-recalculation update
# across sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR (GSM unchanged).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2222.5
$ws.Range("I40").Value = 1320
$ws.Range("J40").Value = 3125
$ws.Range("K40").Value = 1320
$ws.Range("L40").Value = 3125
$ws.Range("M40").Value = -1145
$ws.Range("N40").Value = -3475

$ws.Range("H70").Value = 1440
$ws.Range("I70").Value = 1566.6666
$ws.Range("K70").Value = 4699.9998
$ws.Range("M70").Value = -4429.9998

$ws.Range("H73").Value = 1440
$ws.Range("I73").Value = 1566.6666
$ws.Range("K73").Value = 4699.9998
$ws.Range("M73").Value = -3763.9998

$ws.Range("H108").Value = 45684
$ws.Range("J108").Value = 45684
$ws.Range("L108").Value = 45684
$ws.Range("N108").Value = -53364

$ws.Range("H129").Value = 204930.4
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws.Range("H132").Value = 3189.9656
$ws.Range("I132").Value = 3220.12
$ws.Range("K132").Value = 9660.360000000001
$ws.Range("M132").Value = -7130.360000000001

$ws.Range("H138").Value = 1628.1974
$ws.Range("I138").Value = 1140.6111
$ws.Range("J138").Value = 2067.025
$ws.Range("K138").Value = 3421.8333
$ws.Range("L138").Value = 6201.075000000001
$ws.Range("M138").Value = 1718.1667
$ws.Range("N138").Value = -16481.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16903.984
$ws.Range("I32").Value = 17153.277
$ws.Range("K32").Value = 17153.277
$ws.Range("M32").Value = -16866.277

$ws.Range("H35").Value = 4500
$ws.Range("I35").Value = 4500
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4500
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4094
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1662.6222
$ws.Range("I86").Value = 1411.4062
$ws.Range("K86").Value = 1411.4062
$ws.Range("M86").Value = -288.4061999999999

$ws.Range("H89").Value = 1662.6222
$ws.Range("I89").Value = 1411.4062
$ws.Range("K89").Value = 7057.030999999999
$ws.Range("M89").Value = -1441.030999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12560.37
$ws.Range("I31").Value = 14144.782
$ws.Range("J31").Value = 3450
$ws.Range("K31").Value = 14144.782
$ws.Range("L31").Value = 3450
$ws.Range("M31").Value = -13849.782
$ws.Range("N31").Value = -4040

$ws.Range("H34").Value = 12560.37
$ws.Range("I34").Value = 14144.782
$ws.Range("J34").Value = 3450
$ws.Range("K34").Value = 14144.782
$ws.Range("L34").Value = 3450
$ws.Range("M34").Value = -13942.782
$ws.Range("N34").Value = -3854

$ws.Range("H62").Value = 71431930
$ws.Range("I62").Value = 166669440
$ws.Range("J62").Value = 3795.25
$ws.Range("K62").Value = 166669440
$ws.Range("L62").Value = 3795.25
$ws.Range("M62").Value = -166668816
$ws.Range("N62").Value = -5043.25

$ws.Range("H65").Value = 71431930
$ws.Range("I65").Value = 166669440
$ws.Range("J65").Value = 3795.25
$ws.Range("K65").Value = 833347200
$ws.Range("L65").Value = 18976.25
$ws.Range("M65").Value = -833344080
$ws.Range("N65").Value = -25216.25

$ws.Range("H68").Value = 61154
$ws.Range("J68").Value = 61154
$ws.Range("L68").Value = 61154
$ws.Range("N68").Value = -62652

$ws.Range("H71").Value = 61154
$ws.Range("J71").Value = 61154
$ws.Range("L71").Value = 183462
$ws.Range("N71").Value = -190950

$ws.Range("H86").Value = 5562605.5
$ws.Range("I86").Value = 2349.0625
$ws.Range("K86").Value = 2349.0625
$ws.Range("M86").Value = -1226.0625

$ws.Range("H89").Value = 5562605.5
$ws.Range("I89").Value = 2349.0625
$ws.Range("K89").Value = 11745.3125
$ws.Range("M89").Value = -6129.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2407.8333
$ws.Range("I2").Value = 2881.743
$ws.Range("J2").Value = 38.285713
$ws.Range("K2").Value = 17290.458
$ws.Range("L2").Value = 229.714278
$ws.Range("M2").Value = -17177.458
$ws.Range("N2").Value = -455.714278

$ws.Range("H26").Value = 610.8570999999999
$ws.Range("J26").Value = 763.2
$ws.Range("L26").Value = 2289.6
$ws.Range("N26").Value = -2865.6

$ws.Range("H47").Value = 646
$ws.Range("I47").Value = 76.666664
$ws.Range("J47").Value = 1500
$ws.Range("K47").Value = 229.999992
$ws.Range("L47").Value = 4500
$ws.Range("M47").Value = 201.000008
$ws.Range("N47").Value = -5362

$ws.Range("H109").Value = 1407
$ws.Range("I109").Value = 459.33334
$ws.Range("K109").Value = 1378.00002
$ws.Range("M109").Value = -338.0000199999999

$ws.Range("H110").Value = 2371.8
$ws.Range("I110").Value = 286.66666
$ws.Range("J110").Value = 5499.5
$ws.Range("K110").Value = 859.9999799999999
$ws.Range("L110").Value = 16498.5
$ws.Range("M110").Value = 3230.00002
$ws.Range("N110").Value = -24678.5

$ws.Range("H112").Value = 3266.5
$ws.Range("J112").Value = 3266.5
$ws.Range("L112").Value = 9799.5
$ws.Range("N112").Value = -12015.5

$ws.Range("H114").Value = 967.94116
$ws.Range("J114").Value = 713.2222
$ws.Range("L114").Value = 2139.6666
$ws.Range("N114").Value = -8647.6666

$ws.Range("H118").Value = 166668050
$ws.Range("I118").Value = 250000060
$ws.Range("J118").Value = 4000
$ws.Range("K118").Value = 750000180
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = -749998937
$ws.Range("N118").Value = -14486

$ws.Range("H119").Value = 4383.3335
$ws.Range("I119").Value = 3260
$ws.Range("K119").Value = 9780
$ws.Range("M119").Value = -4942

$ws.Range("H120").Value = 11173.75

$ws.Range("H121").Value = 4115.567
$ws.Range("I121").Value = 441
$ws.Range("J121").Value = 4850.48
$ws.Range("K121").Value = 1323
$ws.Range("L121").Value = 14551.44
$ws.Range("M121").Value = -13
$ws.Range("N121").Value = -17171.44

$ws.Range("H131").Value = 753.5599999999999
$ws.Range("J131").Value = 753.5599999999999
$ws.Range("L131").Value = 2260.68
$ws.Range("N131").Value = -12340.68

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2783.8845
$ws.Range("I68").Value = 2730.3125
$ws.Range("J68").Value = 2869.6
$ws.Range("K68").Value = 2730.3125
$ws.Range("L68").Value = 2869.6
$ws.Range("M68").Value = -1981.3125
$ws.Range("N68").Value = -4367.6

$ws.Range("H71").Value = 2783.8845
$ws.Range("I71").Value = 2730.3125
$ws.Range("J71").Value = 2869.6
$ws.Range("K71").Value = 13651.5625
$ws.Range("L71").Value = 14348
$ws.Range("M71").Value = -9907.5625
$ws.Range("N71").Value = -21836

$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1388.6666
$ws.Range("I122").Value = 1412.1666
$ws.Range("J122").Value = 1329.9166
$ws.Range("K122").Value = 4236.4998
$ws.Range("L122").Value = 3989.7498
$ws.Range("M122").Value = -1786.4998
$ws.Range("N122").Value = -8889.7498
